# Update countries & provincias Spain
# Applies the COVID-19 stats refresh captured in the commit diff:
#  - Updates the "Datos actualizados a ..." timestamp string (A1).
#  - Refreshes Casos totales / Nuevos casos / Casos activos / Recuperados /
#    Casos criticos / Muertes hoy / Muertes for the rows whose data changed.
#  - A handful of countries swapped row order (the new-data row moved above
#    its neighbour), so both rows of each pair are rewritten (name + data)
#    to land on the correct row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: refreshed timestamp ---------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 19 de Septiembre de 2020 a las 10:36"

# --- Row 7: Rusia (data refresh only) ----------------------------------
$ws.Cells.Item(7, 1).Value = "Rusia"
$ws.Cells.Item(7, 2).Value = 1097251
$ws.Cells.Item(7, 3).Value = 6065
$ws.Cells.Item(7, 4).Value = 906462
$ws.Cells.Item(7, 5).Value = 171450
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 144
$ws.Cells.Item(7, 8).Value = 19339

# --- Row 24: Filipinas (data refresh only) -----------------------------
$ws.Cells.Item(24, 1).Value = "Filipinas"
$ws.Cells.Item(24, 2).Value = 283460
$ws.Cells.Item(24, 3).Value = 3962
$ws.Cells.Item(24, 4).Value = 209885
$ws.Cells.Item(24, 5).Value = 68645
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 100
$ws.Cells.Item(24, 8).Value = 4930

# --- Row 26: Indonesia (data refresh only) -----------------------------
$ws.Cells.Item(26, 1).Value = "Indonesia"
$ws.Cells.Item(26, 2).Value = 240687
$ws.Cells.Item(26, 3).Value = 4168
$ws.Cells.Item(26, 4).Value = 174350
$ws.Cells.Item(26, 5).Value = 56889
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 112
$ws.Cells.Item(26, 8).Value = 9448

# --- Row 57: Singapur (data refresh only) ------------------------------
$ws.Cells.Item(57, 1).Value = "Singapur"
$ws.Cells.Item(57, 2).Value = 57558
$ws.Cells.Item(57, 3).Value = 15
$ws.Cells.Item(57, 4).Value = 57071
$ws.Cells.Item(57, 5).Value = 460
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 0
$ws.Cells.Item(57, 8).Value = 27

# --- Rows 69/70: Austria moves above Kenia -----------------------------
$ws.Cells.Item(69, 1).Value = "Austria"
$ws.Cells.Item(69, 2).Value = 37474
$ws.Cells.Item(69, 3).Value = 813
$ws.Cells.Item(69, 4).Value = 28961
$ws.Cells.Item(69, 5).Value = 7748
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 8).Value = 765

$ws.Cells.Item(70, 1).Value = "Kenia"
$ws.Cells.Item(70, 2).Value = 36724
$ws.Cells.Item(70, 3).Value = 0
$ws.Cells.Item(70, 4).Value = 23709
$ws.Cells.Item(70, 5).Value = 12369
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 646

# --- Row 76: El Salvador (data refresh only) ---------------------------
$ws.Cells.Item(76, 1).Value = "El Salvador"
$ws.Cells.Item(76, 2).Value = 27428
$ws.Cells.Item(76, 3).Value = 82
$ws.Cells.Item(76, 4).Value = 21247
$ws.Cells.Item(76, 5).Value = 5373
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 4
$ws.Cells.Item(76, 8).Value = 808

# --- Rows 89/90: Croacia moves above Senegal ---------------------------
$ws.Cells.Item(89, 1).Value = "Croacia"
$ws.Cells.Item(89, 2).Value = 14725
$ws.Cells.Item(89, 3).Value = 212
$ws.Cells.Item(89, 4).Value = 12353
$ws.Cells.Item(89, 5).Value = 2128
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 244

$ws.Cells.Item(90, 1).Value = "Senegal"
$ws.Cells.Item(90, 2).Value = 14645
$ws.Cells.Item(90, 3).Value = 0
$ws.Cells.Item(90, 4).Value = 11051
$ws.Cells.Item(90, 5).Value = 3293
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 301

# --- Row 144: Estonia (data refresh only) -------------------------------
$ws.Cells.Item(144, 1).Value = "Estonia"
$ws.Cells.Item(144, 2).Value = 2875
$ws.Cells.Item(144, 3).Value = 61
$ws.Cells.Item(144, 4).Value = 2374
$ws.Cells.Item(144, 5).Value = 437
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 64

# --- Row 161: Letonia (data refresh only) -------------------------------
$ws.Cells.Item(161, 1).Value = "Letonia"
$ws.Cells.Item(161, 2).Value = 1515
$ws.Cells.Item(161, 3).Value = 17
$ws.Cells.Item(161, 4).Value = 1248
$ws.Cells.Item(161, 5).Value = 231
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 36

# --- Rows 204/205: Santa Lucia moves above Timor Oriental ---------------
$ws.Cells.Item(204, 1).Value = "Santa Lucia"
$ws.Cells.Item(204, 2).Value = 27
$ws.Cells.Item(204, 3).Value = 0
$ws.Cells.Item(204, 4).Value = 26
$ws.Cells.Item(204, 5).Value = 1
$ws.Cells.Item(204, 6).Value = 0
$ws.Cells.Item(204, 7).Value = 0
$ws.Cells.Item(204, 8).Value = 0

$ws.Cells.Item(205, 1).Value = "Timor Oriental"
$ws.Cells.Item(205, 2).Value = 27
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 26
$ws.Cells.Item(205, 5).Value = 1
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 0

# --- Rows 214/215: Montserrat moves above Islas Malvinas ----------------
$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(214, 2).Value = 13
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 1

$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 2).Value = 13
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0
